$d = $word.ActiveDocument

# Each call replaces one exact "a×b=c" answer cell with its corrected value.
# MatchCase=true, MatchWholeWord=false, MatchWildcards=false,
# MatchSoundsLike=false, MatchAllWordForms=false, Forward=true,
# Wrap=wdFindContinue(1), Format=false, Replace=wdReplaceOne(2)

$ok = $d.Content.Find.Execute("71×54=3834", $true, $false, $false, $false, $false, $true, 1, $false, "37×21=777", 2)
if (-not $ok) { throw "Replace failed: 71×54=3834" }
Write-Host "Replaced 71×54=3834 -> 37×21=777: $ok"
$ok = $d.Content.Find.Execute("77×46=3542", $true, $false, $false, $false, $false, $true, 1, $false, "71×95=6745", 2)
if (-not $ok) { throw "Replace failed: 77×46=3542" }
Write-Host "Replaced 77×46=3542 -> 71×95=6745: $ok"
$ok = $d.Content.Find.Execute("76×88=6688", $true, $false, $false, $false, $false, $true, 1, $false, "98×98=9604", 2)
if (-not $ok) { throw "Replace failed: 76×88=6688" }
Write-Host "Replaced 76×88=6688 -> 98×98=9604: $ok"
$ok = $d.Content.Find.Execute("96×38=3648", $true, $false, $false, $false, $false, $true, 1, $false, "89×25=2225", 2)
if (-not $ok) { throw "Replace failed: 96×38=3648" }
Write-Host "Replaced 96×38=3648 -> 89×25=2225: $ok"
$ok = $d.Content.Find.Execute("30×64=1920", $true, $false, $false, $false, $false, $true, 1, $false, "93×34=3162", 2)
if (-not $ok) { throw "Replace failed: 30×64=1920" }
Write-Host "Replaced 30×64=1920 -> 93×34=3162: $ok"
$ok = $d.Content.Find.Execute("58×31=1798", $true, $false, $false, $false, $false, $true, 1, $false, "32×76=2432", 2)
if (-not $ok) { throw "Replace failed: 58×31=1798" }
Write-Host "Replaced 58×31=1798 -> 32×76=2432: $ok"
$ok = $d.Content.Find.Execute("66×32=2112", $true, $false, $false, $false, $false, $true, 1, $false, "27×20=540", 2)
if (-not $ok) { throw "Replace failed: 66×32=2112" }
Write-Host "Replaced 66×32=2112 -> 27×20=540: $ok"
$ok = $d.Content.Find.Execute("12×88=1056", $true, $false, $false, $false, $false, $true, 1, $false, "57×74=4218", 2)
if (-not $ok) { throw "Replace failed: 12×88=1056" }
Write-Host "Replaced 12×88=1056 -> 57×74=4218: $ok"
$ok = $d.Content.Find.Execute("54×91=4914", $true, $false, $false, $false, $false, $true, 1, $false, "91×12=1092", 2)
if (-not $ok) { throw "Replace failed: 54×91=4914" }
Write-Host "Replaced 54×91=4914 -> 91×12=1092: $ok"
$ok = $d.Content.Find.Execute("98×75=7350", $true, $false, $false, $false, $false, $true, 1, $false, "61×19=1159", 2)
if (-not $ok) { throw "Replace failed: 98×75=7350" }
Write-Host "Replaced 98×75=7350 -> 61×19=1159: $ok"
$ok = $d.Content.Find.Execute("25×17=425", $true, $false, $false, $false, $false, $true, 1, $false, "60×62=3720", 2)
if (-not $ok) { throw "Replace failed: 25×17=425" }
Write-Host "Replaced 25×17=425 -> 60×62=3720: $ok"
$ok = $d.Content.Find.Execute("63×99=6237", $true, $false, $false, $false, $false, $true, 1, $false, "88×46=4048", 2)
if (-not $ok) { throw "Replace failed: 63×99=6237" }
Write-Host "Replaced 63×99=6237 -> 88×46=4048: $ok"
$ok = $d.Content.Find.Execute("64×60=3840", $true, $false, $false, $false, $false, $true, 1, $false, "36×78=2808", 2)
if (-not $ok) { throw "Replace failed: 64×60=3840" }
Write-Host "Replaced 64×60=3840 -> 36×78=2808: $ok"
$ok = $d.Content.Find.Execute("92×54=4968", $true, $false, $false, $false, $false, $true, 1, $false, "73×87=6351", 2)
if (-not $ok) { throw "Replace failed: 92×54=4968" }
Write-Host "Replaced 92×54=4968 -> 73×87=6351: $ok"
$ok = $d.Content.Find.Execute("23×87=2001", $true, $false, $false, $false, $false, $true, 1, $false, "71×94=6674", 2)
if (-not $ok) { throw "Replace failed: 23×87=2001" }
Write-Host "Replaced 23×87=2001 -> 71×94=6674: $ok"
$ok = $d.Content.Find.Execute("35×44=1540", $true, $false, $false, $false, $false, $true, 1, $false, "30×31=930", 2)
if (-not $ok) { throw "Replace failed: 35×44=1540" }
Write-Host "Replaced 35×44=1540 -> 30×31=930: $ok"
$ok = $d.Content.Find.Execute("35×90=3150", $true, $false, $false, $false, $false, $true, 1, $false, "49×86=4214", 2)
if (-not $ok) { throw "Replace failed: 35×90=3150" }
Write-Host "Replaced 35×90=3150 -> 49×86=4214: $ok"
$ok = $d.Content.Find.Execute("28×61=1708", $true, $false, $false, $false, $false, $true, 1, $false, "93×18=1674", 2)
if (-not $ok) { throw "Replace failed: 28×61=1708" }
Write-Host "Replaced 28×61=1708 -> 93×18=1674: $ok"
$ok = $d.Content.Find.Execute("85×26=2210", $true, $false, $false, $false, $false, $true, 1, $false, "89×34=3026", 2)
if (-not $ok) { throw "Replace failed: 85×26=2210" }
Write-Host "Replaced 85×26=2210 -> 89×34=3026: $ok"
$ok = $d.Content.Find.Execute("41×97=3977", $true, $false, $false, $false, $false, $true, 1, $false, "93×65=6045", 2)
if (-not $ok) { throw "Replace failed: 41×97=3977" }
Write-Host "Replaced 41×97=3977 -> 93×65=6045: $ok"
$ok = $d.Content.Find.Execute("58×69=4002", $true, $false, $false, $false, $false, $true, 1, $false, "56×67=3752", 2)
if (-not $ok) { throw "Replace failed: 58×69=4002" }
Write-Host "Replaced 58×69=4002 -> 56×67=3752: $ok"
$ok = $d.Content.Find.Execute("39×35=1365", $true, $false, $false, $false, $false, $true, 1, $false, "51×98=4998", 2)
if (-not $ok) { throw "Replace failed: 39×35=1365" }
Write-Host "Replaced 39×35=1365 -> 51×98=4998: $ok"
$ok = $d.Content.Find.Execute("77×81=6237", $true, $false, $false, $false, $false, $true, 1, $false, "35×91=3185", 2)
if (-not $ok) { throw "Replace failed: 77×81=6237" }
Write-Host "Replaced 77×81=6237 -> 35×91=3185: $ok"
$ok = $d.Content.Find.Execute("39×29=1131", $true, $false, $false, $false, $false, $true, 1, $false, "75×59=4425", 2)
if (-not $ok) { throw "Replace failed: 39×29=1131" }
Write-Host "Replaced 39×29=1131 -> 75×59=4425: $ok"
$ok = $d.Content.Find.Execute("96×46=4416", $true, $false, $false, $false, $false, $true, 1, $false, "13×68=884", 2)
if (-not $ok) { throw "Replace failed: 96×46=4416" }
Write-Host "Replaced 96×46=4416 -> 13×68=884: $ok"
